$wb = $excel.ActiveWorkbook

# --- Step 1: capture the OLD "dept hours" data before we overwrite it ---
$deptSheet = $wb.Worksheets.Item("dept hours")

$oldRows = @()
for ($r = 2; $r -le 5; $r++) {
    $b = $deptSheet.Cells.Item($r, 2).Value2
    $c = $deptSheet.Cells.Item($r, 3).Value2
    $d = $deptSheet.Cells.Item($r, 4).Value2
    $oldRows += ,@($b, $c, $d)
}

# --- Step 2: update "PI hours" sheet - add "app" column (old dept lists), and set dept to single value ---
$piSheet = $wb.Worksheets.Item("PI hours")

$oldE2 = $piSheet.Range("E2").Value2
$oldE3 = $piSheet.Range("E3").Value2

# Give F1 the same header style/formatting as E1
$piSheet.Range("E1").Copy()
$piSheet.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$piSheet.Range("F1").Value2 = "app"

# Move the old list values from E into F, then set new single-dept values into E
$piSheet.Range("F2").Value2 = $oldE2
$piSheet.Range("F3").Value2 = $oldE3

$piSheet.Range("E2").Value2 = "ME"
$piSheet.Range("E3").Value2 = "ECE"

# --- Step 3: rename "dept hours" to "department hours" and replace its content ---
$deptSheet.Name = "department hours"

# Delete old rows 4 and 5 (sheet shrinks from 5 rows to 3 rows)
$deptSheet.Rows.Item(5).Delete()
$deptSheet.Rows.Item(4).Delete()

# New simplified content: one row per dept that actually appears as a PI's dept
$deptSheet.Range("B2").Value2 = "ME"
$deptSheet.Range("C2").Value2 = 63.5
$deptSheet.Range("D2").Value2 = 63.18407960199005

$deptSheet.Range("B3").Value2 = "ECE"
$deptSheet.Range("C3").Value2 = 37
$deptSheet.Range("D3").Value2 = 36.81592039800995

# --- Step 4: add the new "unit(accumulative) hours" sheet with the OLD dept-hours data ---
$unitSheet = $wb.Worksheets.Add()
$unitSheet.Name = "unit(accumulative) hours"

$deptSheet.Range("B1:D1").Copy()
$unitSheet.Range("B1:D1").PasteSpecial(-4122)  # xlPasteFormats
$unitSheet.Range("B1").Value2 = "unit(accumulative)"
$unitSheet.Range("C1").Value2 = "hours"
$unitSheet.Range("D1").Value2 = "percentage"

for ($i = 0; $i -lt $oldRows.Count; $i++) {
    $row = $i + 2
    $deptSheet.Range("A2").Copy()
    $unitSheet.Cells.Item($row, 1).PasteSpecial(-4122)  # xlPasteFormats
    $unitSheet.Cells.Item($row, 1).Value2 = $i
    $unitSheet.Cells.Item($row, 2).Value2 = $oldRows[$i][0]
    $unitSheet.Cells.Item($row, 3).Value2 = $oldRows[$i][1]
    $unitSheet.Cells.Item($row, 4).Value2 = $oldRows[$i][2]
}

# Move the new sheet to the end (after "department hours")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$unitSheet.Move($null, $lastSheet)
